$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 21
$ws.Range("F4").Value = 7780
$ws.Range("F5").Value = 99
$ws.Range("F6").Value = 117
$ws.Range("F8").Value = 8585
$ws.Range("F9").Value = 8585
$ws.Range("F10").Value = 12
$ws.Range("F12").Value = 88
$ws.Range("F13").Value = 5728
$ws.Range("F14").Value = 62
$ws.Range("F15").Value = 2687
$ws.Range("F16").Value = 1187
$ws.Range("F20").Value = 39
$ws.Range("F21").Value = 581
$ws.Range("F22").Value = 12
$ws.Range("F23").Value = 3755
$ws.Range("F24").Value = 75
$ws.Range("F25").Value = 50
$ws.Range("F26").Value = 42
$ws.Range("F27").Value = 13
$ws.Range("F28").Value = 6
$ws.Range("F29").Value = 3481
$ws.Range("F33").Value = 374
$ws.Range("F34").Value = 150
$ws.Range("F35").Value = 359
$ws.Range("F36").Value = 1299
$ws.Range("F37").Value = 690
$ws.Range("F40").Value = 3067
$ws.Range("F41").Value = 59
$ws.Range("F44").Value = 3323
$ws.Range("F46").Value = 2313

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 115
$ws.Range("F3").Value = 141
$ws.Range("F5").Value = 61

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1353

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1353
$ws.Range("F4").Value = 21
$ws.Range("F5").Value = 7780
$ws.Range("F6").Value = 99
$ws.Range("F7").Value = 117
$ws.Range("F9").Value = 8585
$ws.Range("F10").Value = 8585
$ws.Range("F11").Value = 12
$ws.Range("F12").Value = 88
$ws.Range("F13").Value = 5728
$ws.Range("F14").Value = 62
$ws.Range("F15").Value = 2687
$ws.Range("F16").Value = 1187
$ws.Range("F20").Value = 115
$ws.Range("F21").Value = 39
$ws.Range("F22").Value = 141
$ws.Range("F23").Value = 581
$ws.Range("F25").Value = 3755
$ws.Range("F26").Value = 75
$ws.Range("F27").Value = 50
$ws.Range("F28").Value = 42
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 3483
$ws.Range("F32").Value = 374
$ws.Range("F33").Value = 150
$ws.Range("F34").Value = 359
$ws.Range("F35").Value = 61
$ws.Range("F36").Value = 1300
$ws.Range("F37").Value = 690
$ws.Range("F41").Value = 3067
$ws.Range("F42").Value = 59
$ws.Range("F45").Value = 3323
$ws.Range("F47").Value = 2313
